$d = $word.ActiveDocument

# The final paragraph currently reads "Manejo de Glasfish" and carries the
# "_GoBack" bookmark. We need to:
#   1. Split it into two runs ("Manejo de " / "Glasfish") wrapped with
#      spell-check proofErr markers (no bookmark on this paragraph anymore).
#   2. Add a new, empty paragraph after it.
#   3. Add a further paragraph with the new explanatory text (including
#      proofErr spans around the lowercase "glasfish" and the ellipsis), and
#      move the "_GoBack" bookmark to the very end of that paragraph's text.

$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# Drop the existing _GoBack bookmark - it will be re-created in its new
# location further down.
try { $d.Bookmarks("_GoBack").Delete() } catch { }

# Remove the paragraph's text content but keep the paragraph mark itself so
# the document structure (and the final section properties) stays intact.
$textOnly = $d.Range($lastRange.Start, $lastRange.End - 1)
$textOnly.Delete()

# Insert point: start of the now-empty trailing paragraph.
$insertAt = $d.Paragraphs.Last.Range.Start
$target = $d.Range($insertAt, $insertAt)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newXml = "<w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">Manejo de </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Glasfish</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
  "</w:p>" +
  "<w:p $wns/>" +
  "<w:p $wns>" +
    "<w:r><w:t xml:space=`"preserve`">Este es un manual para la configuración de </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:proofErr w:type=`"gramStart`"/>" +
    "<w:r><w:t>glasfish</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> …</w:t></w:r>" +
    "<w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r><w:t>.GENERADO DESDE MI MAQUINA LOCAL</w:t></w:r>" +
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
    "<w:bookmarkEnd w:id=`"0`"/>" +
  "</w:p>"

$target.InsertXML($newXml)
